$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the years to keep (2010-2012), which will become rows 2-4
$years  = @("2010年", "2011年", "2012年")
$colB   = @(5.6, 8.3032826787889, 8.167441309779999)
$colC   = @(12515, 22642.985251342, 17796.176422466)
$colD   = @(7219361.8, 10074665.7, 8146746.1)
$colE   = @(3207.1874, 3694.4244, 3738.897)
$colF   = @(2244827, 1246856, 1260582)

for ($i = 0; $i -lt $years.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $years[$i]
    $ws.Cells.Item($r, 2).Value = $colB[$i]
    $ws.Cells.Item($r, 3).Value = $colC[$i]
    $ws.Cells.Item($r, 4).Value = $colD[$i]
    $ws.Cells.Item($r, 5).Value = $colE[$i]
    $ws.Cells.Item($r, 6).Value = $colF[$i]
}

# Remove the now-obsolete rows 5 through 14 (old 2003年-2012年 rows)
$ws.Range("A5:F14").EntireRow.Delete()
